# Adds two new BOM rows (49 & 50) to the "BOM" sheet:
#   Row 49: Waste water tap (RMTL)
#   Row 50: Plumbing bracket mounting screws (Sourcingmap)
# and updates the view/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# ---------------------------------------------------------------------------
# Row 49 - Waste water tap
# ---------------------------------------------------------------------------

# Copy cell formatting down from row 48 (the previous last row) before we
# write any values, so the new cells pick up the same styles (wrap text on
# D, hyperlink style on F, currency format on J/K) without creating
# duplicate style entries.
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("J48:K48").Copy() | Out-Null
$ws.Range("J49:K49").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B49").Value = "Waste water tap"
$ws.Range("C49").Value = "RMTL"
$ws.Range("D49").Value = "RMTL Outdoor Garden Tap Hose Union Bib Tap 1/2inch BSP with Compression Wallplate Elbow, Through The Wall 15mm Copper Tube (Single Bib Tap)"

$ws.Range("F49").Value = "Amazon"
$ws.Hyperlinks.Add($ws.Range("F49"), "https://www.amazon.co.uk/dp/B07D5NHKN8") | Out-Null
# Re-apply F48's format so the hyperlink reuses the workbook's existing
# Hyperlink style (xf index) instead of Excel re-deriving a new one.
$ws.Range("F48").Copy() | Out-Null
$ws.Range("F49").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 1
$ws.Range("I49").Formula = "=G49*H49"
$ws.Range("J49").Value = 8.99
$ws.Range("K49").Formula = "=G49*J49"

$ws.Rows.Item(49).RowHeight = 30

# ---------------------------------------------------------------------------
# Row 50 - Plumbing bracket mounting screws
# ---------------------------------------------------------------------------

$ws.Range("D48").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("J48:K48").Copy() | Out-Null
$ws.Range("J50:K50").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C50").Value = "Sourcingmap"
$ws.Range("D50").Value = "Sourcingmap 50 Pcs M2 x 15mm Stainless Steel Phillips Round Head Self Tapping Screws Bolts"
$ws.Range("E50").Value = "a16072200ux0833"

$ws.Range("F50").Value = "Amazon"
$ws.Hyperlinks.Add($ws.Range("F50"), "https://www.amazon.co.uk/dp/B01CODNIC2") | Out-Null
$ws.Range("F48").Copy() | Out-Null
$ws.Range("F50").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 50
$ws.Range("I50").Formula = "=G50*H50"
$ws.Range("J50").Value = 5.99
$ws.Range("K50").Formula = "=G50*J50"

# B50 is written last - matches the workbook's shared-string append order.
$ws.Range("B50").Value = "Plumbing bracket mounting screws"

# ---------------------------------------------------------------------------
# Clear the leftover copy marquee and update the view/selection so the
# workbook re-opens scrolled to, and selected on, the new last row.
# ---------------------------------------------------------------------------

$excel.CutCopyMode = 0

$ws.Range("B50").Select() | Out-Null
